# Update automatico via Actualizar 02-17-2021 13-10-09
#
# The sheet keeps a rolling log of "Ultimo" (last-checked) timestamps in
# column D, grouped in blocks of 14 rows (one block per check cycle).
# On each run the newest timestamp is written into the first block
# (rows 2-15) and every older block shifts down to the value that was
# previously held by the block above it.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Newest check: rows 2-15 get "now" (captured at commit time).
for ($r = 2; $r -le 15; $r++) {
    $ws.Cells.Item($r, 4).Value = 44244.54857985241
}

# Previous check shifts down into rows 16-29.
for ($r = 16; $r -le 29; $r++) {
    $ws.Cells.Item($r, 4).Value = 44244.52730865741
}

# Check before that shifts down into rows 30-43.
for ($r = 30; $r -le 43; $r++) {
    $ws.Cells.Item($r, 4).Value = 44244.50605336806
}
